$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "G" = 54.86839566666666; "H" = 164.605187; "I" = 0.6170939026906647; "J" = 0.6170939026906647; "M" = 1.442875; "N" = 4.328625; "O" = 0.02047893724893121; "P" = 0.02047893724893121; "Q" = 79.16823639754166; "R" = 712.514127577875; "S" = 0.01263742730990018; "T" = 0.01263742730990018 }
    3 = @{ "G" = 54.86839566666666; "H" = 164.605187; "I" = 0.6170939026906647; "J" = 0.6170939026906647; "O" = 0.1473796107804731; "P" = 0.1473796107804731; "Q" = 569.7455744220864; "R" = 5127.710169798777; "S" = 0.0909470591935533; "T" = 0.09094705919355328 }
    4 = @{ "G" = 54.86839566666666; "H" = 164.605187; "I" = 0.6170939026906647; "J" = 0.6170939026906647; "M" = 27.934719; "N" = 83.804157; "O" = 0.3964815784233052; "P" = 0.3964815784233051; "Q" = 1532.733214929151; "R" = 13794.59893436236; "S" = 0.2446663645741922; "T" = 0.2446663645741922 }
    5 = @{ "G" = 54.86839566666666; "H" = 164.605187; "I" = 0.6170939026906647; "J" = 0.6170939026906647; "M" = 30.695086; "N" = 92.085258; "O" = 0.4356598735472906; "P" = 0.4356598735472905; "Q" = 1684.190123670361; "R" = 15157.71111303325; "S" = 0.268843051613019; "T" = 0.268843051613019 }
    6 = @{ "I" = 0.06720170646055251; "J" = 0.0672017064605525; "M" = 1.442875; "N" = 4.328625; "O" = 0.02047893724893121; "P" = 0.02047893724893121; "Q" = 8.621444094958333; "R" = 77.592996854625; "S" = 0.00137621952962675; "T" = 0.001376219529626749 }
    7 = @{ "I" = 0.06720170646055251; "J" = 0.0672017064605525; "O" = 0.1473796107804731; "P" = 0.1473796107804731; "S" = 0.009904161341939834; "T" = 0.00990416134193983 }
    8 = @{ "I" = 0.06720170646055251; "J" = 0.0672017064605525; "M" = 27.934719; "N" = 83.804157; "O" = 0.3964815784233052; "P" = 0.3964815784233051; "Q" = 166.915095324869; "R" = 1502.235857923821; "S" = 0.02664423865021948; "T" = 0.02664423865021948 }
    9 = @{ "I" = 0.06720170646055251; "J" = 0.0672017064605525; "M" = 30.695086; "N" = 92.085258; "O" = 0.4356598735472906; "P" = 0.4356598735472905; "Q" = 183.4087969775194; "R" = 1650.679172797674; "S" = 0.02927708693876645; "T" = 0.02927708693876644 }
    10 = @{ "G" = 4.832157666666667; "H" = 14.496473; "I" = 0.05434631351453007; "J" = 0.05434631351453007; "M" = 1.442875; "N" = 4.328625; "O" = 0.02047893724893121; "P" = 0.02047893724893121; "Q" = 6.972199493291667; "R" = 62.749795439625; "S" = 0.001112954744174803; "T" = 0.001112954744174803 }
    11 = @{ "G" = 4.832157666666667; "H" = 14.496473; "I" = 0.05434631351453007; "J" = 0.05434631351453007; "O" = 0.1473796107804731; "P" = 0.1473796107804731; "Q" = 50.17643421212034; "R" = 451.5879079090831; "S" = 0.008009538533125008; "T" = 0.008009538533125006 }
    12 = @{ "G" = 4.832157666666667; "H" = 14.496473; "I" = 0.05434631351453007; "J" = 0.05434631351453007; "M" = 27.934719; "N" = 83.804157; "O" = 0.3964815784233052; "P" = 0.3964815784233051; "Q" = 134.984966582029; "R" = 1214.864699238261; "S" = 0.02154731216372869; "T" = 0.02154731216372868 }
    13 = @{ "G" = 4.832157666666667; "H" = 14.496473; "I" = 0.05434631351453007; "J" = 0.05434631351453007; "M" = 30.695086; "N" = 92.085258; "O" = 0.4356598735472906; "P" = 0.4356598735472905; "Q" = 148.3234951438927; "R" = 1334.911456295034; "S" = 0.02367650807350158; "T" = 0.02367650807350158 }
    14 = @{ "G" = 23.238438; "H" = 69.715314; "I" = 0.2613580773342528; "J" = 0.2613580773342528; "M" = 1.442875; "N" = 4.328625; "O" = 0.02047893724893121; "P" = 0.02047893724893121; "Q" = 33.53016122925; "R" = 301.77145106325; "S" = 0.005352335665229472; "T" = 0.005352335665229472 }
    15 = @{ "G" = 23.238438; "H" = 69.715314; "I" = 0.2613580773342528; "J" = 0.2613580773342528; "O" = 0.1473796107804731; "P" = 0.1473796107804731; "Q" = 241.304617095366; "R" = 2171.741553858294; "S" = 0.03851885171185496; "T" = 0.03851885171185495 }
    16 = @{ "G" = 23.238438; "H" = 69.715314; "I" = 0.2613580773342528; "J" = 0.2613580773342528; "M" = 27.934719; "N" = 83.804157; "O" = 0.3964815784233052; "P" = 0.3964815784233051; "Q" = 649.1592355289221; "R" = 5842.433119760299; "S" = 0.1036236630351648; "T" = 0.1036236630351648 }
    17 = @{ "G" = 23.238438; "H" = 69.715314; "I" = 0.2613580773342528; "J" = 0.2613580773342528; "M" = 30.695086; "N" = 92.085258; "O" = 0.4356598735472906; "P" = 0.4356598735472905; "Q" = 713.3058529156681; "R" = 6419.752676241013; "S" = 0.1138632269220036; "T" = 0.1138632269220035 }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}

Write-Output "Updated $($data.Count) rows"